$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 8003.1665
$ws.Range("J48").Value = 8003.1665
$ws.Range("L48").Value = 24009.4995
$ws.Range("N48").Value = -24593.4995
$ws.Range("H56").Value = 8003.1665
$ws.Range("J56").Value = 8003.1665
$ws.Range("L56").Value = 24009.4995
$ws.Range("N56").Value = -25077.4995
$ws.Range("H62").Value = 44184.36
$ws.Range("I62").Value = 86507.414
$ws.Range("J62").Value = 5116.923
$ws.Range("K62").Value = 86507.414
$ws.Range("L62").Value = 5116.923
$ws.Range("M62").Value = -85883.414
$ws.Range("N62").Value = -6364.923
$ws.Range("H65").Value = 44184.36
$ws.Range("I65").Value = 86507.414
$ws.Range("J65").Value = 5116.923
$ws.Range("K65").Value = 432537.07
$ws.Range("L65").Value = 25584.615
$ws.Range("M65").Value = -429417.07
$ws.Range("N65").Value = -31824.615
$ws.Range("H118").Value = 2440.7144
$ws.Range("I118").Value = 4394
$ws.Range("J118").Value = 1355.5555
$ws.Range("K118").Value = 13182
$ws.Range("L118").Value = 4066.6665
$ws.Range("M118").Value = -11525
$ws.Range("N118").Value = -7380.666499999999
$ws.Range("H123").Value = 43000
$ws.Range("J123").Value = 43000
$ws.Range("L123").Value = 43000
$ws.Range("N123").Value = -52800
$ws.Range("H128").Value = 34537.145
$ws.Range("J128").Value = 34537.145
$ws.Range("L128").Value = 34537.145
$ws.Range("N128").Value = -44497.145
$ws.Range("H132").Value = 1572.2456
$ws.Range("I132").Value = 1622.8163
$ws.Range("J132").Value = 1262.5
$ws.Range("K132").Value = 4868.448899999999
$ws.Range("L132").Value = 3787.5
$ws.Range("M132").Value = -2338.448899999999
$ws.Range("N132").Value = -8847.5
$ws.Range("H133").Value = 43353.848
$ws.Range("J133").Value = 43353.848
$ws.Range("L133").Value = 43353.848
$ws.Range("N133").Value = -53473.848
$ws.Range("H136").Value = 41097.95
$ws.Range("I136").Value = 67000
$ws.Range("J136").Value = 39734.684
$ws.Range("K136").Value = 67000
$ws.Range("L136").Value = 39734.684
$ws.Range("M136").Value = -61900
$ws.Range("N136").Value = -49934.684
$ws.Range("H139").Value = 62540
$ws.Range("J139").Value = 62540
$ws.Range("L139").Value = 62540
$ws.Range("N139").Value = -72820
$ws.Range("H141").Value = 2212.4888
$ws.Range("I141").Value = 773.6667
$ws.Range("J141").Value = 7967.778
$ws.Range("K141").Value = 2321.0001
$ws.Range("L141").Value = 23903.334
$ws.Range("M141").Value = 2858.9999
$ws.Range("N141").Value = -34263.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1040.7037
$ws.Range("I61").Value = 847.8421
$ws.Range("K61").Value = 847.8421
$ws.Range("M61").Value = -635.8421
$ws.Range("H74").Value = 986.4146
$ws.Range("I74").Value = 959.0263
$ws.Range("J74").Value = 1333.3334
$ws.Range("K74").Value = 959.0263
$ws.Range("L74").Value = 1333.3334
$ws.Range("M74").Value = -85.02629999999999
$ws.Range("N74").Value = -3081.3334
$ws.Range("H77").Value = 986.4146
$ws.Range("I77").Value = 959.0263
$ws.Range("J77").Value = 1333.3334
$ws.Range("K77").Value = 4795.1315
$ws.Range("L77").Value = 6666.666999999999
$ws.Range("M77").Value = -427.1314999999995
$ws.Range("N77").Value = -15402.667
$ws.Range("H132").Value = 1662.25
$ws.Range("I132").Value = 1471.625
$ws.Range("J132").Value = 2424.75
$ws.Range("K132").Value = 4414.875
$ws.Range("L132").Value = 7274.25
$ws.Range("M132").Value = -1884.875
$ws.Range("N132").Value = -12334.25
$ws.Range("H136").Value = 1040.7037
$ws.Range("I136").Value = 847.8421
$ws.Range("K136").Value = 2543.5263
$ws.Range("M136").Value = 6.473700000000008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18973.121
$ws.Range("I134").Value = 1538.7609
$ws.Range("J134").Value = 85804.836
$ws.Range("K134").Value = 4616.2827
$ws.Range("L134").Value = 257414.508
$ws.Range("M134").Value = -2081.2827
$ws.Range("N134").Value = -262484.508

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2760.9075
$ws.Range("I58").Value = 1000.2727
$ws.Range("J58").Value = 5527.619
$ws.Range("K58").Value = 1000.2727
$ws.Range("L58").Value = 5527.619
$ws.Range("M58").Value = -797.2727
$ws.Range("N58").Value = -5933.619
$ws.Range("H132").Value = 1888.1404
$ws.Range("I132").Value = 1137.4482
$ws.Range("J132").Value = 2665.6428
$ws.Range("K132").Value = 3412.3446
$ws.Range("L132").Value = 7996.928400000001
$ws.Range("M132").Value = -882.3446000000004
$ws.Range("N132").Value = -13056.9284
$ws.Range("H134").Value = 1485.3914
$ws.Range("I134").Value = 1419.8334
$ws.Range("J134").Value = 1721.4
$ws.Range("K134").Value = 4259.5002
$ws.Range("L134").Value = 5164.200000000001
$ws.Range("M134").Value = -1724.5002
$ws.Range("N134").Value = -10234.2
$ws.Range("H136").Value = 2760.9075
$ws.Range("I136").Value = 1000.2727
$ws.Range("J136").Value = 5527.619
$ws.Range("K136").Value = 3000.8181
$ws.Range("L136").Value = 16582.857
$ws.Range("M136").Value = -450.8181
$ws.Range("N136").Value = -21682.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 596.8372000000001
$ws.Range("I5").Value = 558.62164
$ws.Range("K5").Value = 1675.86492
$ws.Range("M5").Value = -1563.86492
$ws.Range("H135").Value = 596.8372000000001
$ws.Range("I135").Value = 558.62164
$ws.Range("K135").Value = 5027.59476
$ws.Range("M135").Value = -2492.59476

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4048.302
$ws.Range("I70").Value = 3927.1794
$ws.Range("K70").Value = 3927.1794
$ws.Range("M70").Value = -3657.1794
$ws.Range("H73").Value = 4048.302
$ws.Range("I73").Value = 3927.1794
$ws.Range("K73").Value = 3927.1794
$ws.Range("M73").Value = -2991.1794
$ws.Range("H80").Value = 3800
$ws.Range("I80").Value = 3800
$ws.Range("K80").Value = 3800
$ws.Range("M80").Value = -2802
$ws.Range("H82").Value = 30665.834
$ws.Range("J82").Value = 30665.834
$ws.Range("L82").Value = 30665.834
$ws.Range("N82").Value = -31431.834
$ws.Range("H83").Value = 3800
$ws.Range("I83").Value = 3800
$ws.Range("K83").Value = 19000
$ws.Range("M83").Value = -14008
$ws.Range("H85").Value = 30665.834
$ws.Range("J85").Value = 30665.834
$ws.Range("L85").Value = 30665.834
$ws.Range("N85").Value = -33317.834
$ws.Range("H126").Value = 2720.0588
$ws.Range("I126").Value = 2659.0715
$ws.Range("K126").Value = 7977.2145
$ws.Range("M126").Value = -5507.2145
$ws.Range("H132").Value = 2730.3872
$ws.Range("I132").Value = 2579.842
$ws.Range("J132").Value = 2968.75
$ws.Range("K132").Value = 7739.526
$ws.Range("L132").Value = 8906.25
$ws.Range("M132").Value = -5209.526
$ws.Range("N132").Value = -13966.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5805.5
$ws.Range("I16").Value = 8596.866
$ws.Range("J16").Value = 1999.091
$ws.Range("K16").Value = 8596.866
$ws.Range("L16").Value = 1999.091
$ws.Range("M16").Value = -8426.866
$ws.Range("N16").Value = -2339.091
$ws.Range("H132").Value = 2016.0834
$ws.Range("I132").Value = 1738.907
$ws.Range("J132").Value = 4399.8
$ws.Range("K132").Value = 5216.721
$ws.Range("L132").Value = 13199.4
$ws.Range("M132").Value = -2686.721
$ws.Range("N132").Value = -18259.4
$ws.Range("H136").Value = 2157.2727
$ws.Range("I136").Value = 1138.8235
$ws.Range("K136").Value = 3416.4705
$ws.Range("M136").Value = -866.4704999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 729.7646999999999
$ws.Range("I132").Value = 618.08887
$ws.Range("J132").Value = 1567.3334
$ws.Range("K132").Value = 1854.26661
$ws.Range("L132").Value = 4702.0002
$ws.Range("M132").Value = 675.7333899999999
$ws.Range("N132").Value = -9762.0002
$ws.Range("H136").Value = 847.4761999999999
$ws.Range("I136").Value = 961.8214
$ws.Range("J136").Value = 618.7857
$ws.Range("K136").Value = 2885.4642
$ws.Range("L136").Value = 1856.3571
$ws.Range("M136").Value = -335.4642000000003
$ws.Range("N136").Value = -6956.3571
